# Apply the StructureDefinition metadata refresh:
#  - Metadata sheet: URL, Version, Date, Publisher updated (IBM/Alvearie -> LinuxForHealth)
#  - Elements sheet: clear the stray "Constraint(s)" text duplicated on the root Extension row

$wb = $excel.ActiveWorkbook

$wsMetadata = $wb.Worksheets.Item("Metadata")
$wsMetadata.Range("B2").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-exempt-indicator"
$wsMetadata.Range("B3").Value = "8.0.0"
$wsMetadata.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$wsMetadata.Range("B9").Value = "LinuxForHealth Team"

$wsElements = $wb.Worksheets.Item("Elements")
$wsElements.Range("AI2").Value = ""
$wsElements.Range("Q5").Value = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-exempt-indicator"
